$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Scratch cell used to produce a genuine Text-typed value ("04-Nov-2025")
# without Excel's automatic date-recognition kicking in. We build the text
# via a formula (so it's typed as Text), copy it, then PasteSpecial just the
# values onto each target cell - PasteSpecial carries over the source's
# Text type instead of re-parsing the string as a date.
$scratch = $ws.Cells.Item(1, 11)
$scratch.Formula = "=""04-Nov-2025"""
$scratch.Copy()

for ($r = 3; $r -le 19; $r++) {
    $hCell = $ws.Cells.Item($r, 8)
    $hCell.Value = $hCell.Value2 - 1

    $ws.Cells.Item($r, 9).PasteSpecial(-4163)
}

$scratch.Clear()
$excel.CutCopyMode = 0
